# Applies the "homogeneous_evolutionary" solution-80 update:
#  - Resumen!B2/C2: zone + time of the overall best solution change
#  - Solucion!B<row>: reshuffle of the "Salida" (output slot) assignments
#  - Metricas!B2:B5: refreshed per-zone timings

$wb  = $excel.ActiveWorkbook
$wsResumen  = $wb.Worksheets.Item("Resumen")
$wsSolucion = $wb.Worksheets.Item("Solucion")
$wsMetricas = $wb.Worksheets.Item("Metricas")

# --- Resumen sheet: best zone / best time -------------------------------
$wsResumen.Range("B2").Value = "Z4"
$wsResumen.Range("C2").Value = 522.1891231484485

# --- Solucion sheet: re-assign "Salida" values for several Pedido rows --
$solucionUpdates = @{
    11 = "S002"
    13 = "S022"
    15 = "S032"
    16 = "S012"
    19 = "S003"
    20 = "S023"
    21 = "S074"
    23 = "S033"
    25 = "S053"
    28 = "S035"
    30 = "S034"
    31 = "S025"
    32 = "S043"
    33 = "S075"
    34 = "S014"
    36 = "S065"
    37 = "S045"
    38 = "S004"
    40 = "S054"
    42 = "S055"
    43 = "S046"
    45 = "S006"
    50 = "S056"
    52 = "S067"
    54 = "S069"
    55 = "S017"
    56 = "S050"
    61 = "S038"
    63 = "S077"
    66 = "S030"
    67 = "S070"
    74 = "S058"
    75 = "S047"
    77 = "S078"
}

foreach ($row in $solucionUpdates.Keys) {
    $wsSolucion.Range("B$row").Value = $solucionUpdates[$row]
}

# --- Metricas sheet: updated per-zone timings ---------------------------
$wsMetricas.Range("B2").Value = 522.1704941074711
$wsMetricas.Range("B3").Value = 522.1560817385662
$wsMetricas.Range("B4").Value = 522.1607092658669
$wsMetricas.Range("B5").Value = 522.1891231484485
